# Swap the deck's two embedded themes: the slide-master theme (ppt/theme/theme1.xml,
# currently "Integral") becomes the "Office Theme" colour scheme that previously lived
# on the notes master (ppt/theme/theme2.xml).
#
# The PowerPoint object model in this host does not expose a supported way to rewrite
# ppt/theme/theme2.xml (the notes-master theme) directly, nor the <a:theme>/<a:clrScheme>
# "name" attributes - only the 12 theme colour slots are settable, via
# Slide.ThemeColorScheme, and that collection is backed by the slide master's theme part
# (theme1.xml). All slides share one master, so writing through slide 1 updates the
# whole deck's active theme colours.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$tcs = $s.ThemeColorScheme

# Index -> (slot, target "Office Theme" RGB as a COM BGR-packed long: R | (G<<8) | (B<<16))
# 1 dk1       000000 -> 0
# 2 lt1       FFFFFF -> 16777215
# 3 dk2       44546A -> 6968388
# 4 lt2       E7E6E6 -> 15132391
# 5 accent1   5B9BD5 -> 13998939
# 6 accent2   ED7D31 -> 3243501
# 7 accent3   A5A5A5 -> 10855845
# 8 accent4   FFC000 -> 49407
# 9 accent5   4472C4 -> 12874308
# 10 accent6  70AD47 -> 4697456
# 11 hlink    0563C1 -> 12673797
# 12 folHlink 954F72 -> 7491477

$tcs.Colors(1).RGB = 0
$tcs.Colors(2).RGB = 16777215
$tcs.Colors(3).RGB = 6968388
$tcs.Colors(4).RGB = 15132391
$tcs.Colors(5).RGB = 13998939
$tcs.Colors(6).RGB = 3243501
$tcs.Colors(7).RGB = 10855845
$tcs.Colors(8).RGB = 49407
$tcs.Colors(9).RGB = 12874308
$tcs.Colors(10).RGB = 4697456
$tcs.Colors(11).RGB = 12673797
$tcs.Colors(12).RGB = 7491477
